$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: phone number for Amar Harripersad becomes a text value "704-6465"
# instead of a raw number 7046465.
$ws.Range("E21").Value = "704-6465"

# Row 22: Sean-Michael Gopaul gets an email (hyperlink) and phone number.
$ws.Hyperlinks.Add($ws.Range("D22"), "mailto:qertyblue@gmail.com", $null, $null, "qertyblue@gmail.com")
$ws.Range("D22").Style = "Hyperlink"
$ws.Range("E22").Value = "787-4565"

# Row 23: brand new signup - Anirudh Madala.
$ws.Range("A23").Value = "Anirudh Madala"
$ws.Range("B23").Value = "4P"
$ws.Range("C23").Value = "Participant"
$ws.Range("E23").Value = "738-3113"

# Match the final selection/active cell left by the editor.
$ws.Range("D10").Select()
